$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows for the bug-tracking sheet.
# Columns: A=TYPE, B=STATUS, C=DATE, D=SUMMARY
$ws.Range("A4").Value = "bug"
$ws.Range("B4").Value = "open"
$ws.Range("C4").Value = "Webnesday, April 20, 2021"
$ws.Range("D4").Value = "fix scrolling capability in time entry"

$ws.Range("A5").Value = "feature request"
$ws.Range("B5").Value = "open"
$ws.Range("C5").Value = "Webnesday, April 20, 2021"
$ws.Range("D5").Value = "allow CRUD operations on course data"

# Widen column A and move the active selection to match the authored diff.
# (ColumnWidth is quantized on round-trip by the host; 15.65 is the input
# that lands closest to the authored OOXML width of 16.42578125.)
$ws.Columns.Item(1).ColumnWidth = 15.65
$ws.Range("D6").Select()
